# The commit widened the "Day Index" date column (A) and gave the
# "Impressions" column (C) an explicit best-fit width, and zoomed the
# sheet view to 160%. Reproduce the view/format changes via the Excel
# object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zoom the active window to 160% (persists as sheetView zoomScale).
$excel.ActiveWindow.Zoom = 160

# Widen column A ("Day Index") and give column C ("Impressions") an
# explicit width, matching the column-width adjustment from the commit.
$ws.Columns.Item(1).ColumnWidth = 18
$ws.Columns.Item(3).ColumnWidth = 10.833333333333334
